$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.017.58'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '2.256.85'
$ws.Range("E3").Value = '  -1.35%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.522'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.67%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0783'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.76%  '

$ws.Range("E12").Value = '  +0.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.68%  '

$ws.Range("D14").Value = '2.604.52'
$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("D16").Value = '2.264.23'
$ws.Range("E16").Value = '  -0.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.784'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.57%  '

$ws.Range("D18").Value = '41.844.52'
$ws.Range("E18").Value = '  -1.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.02%  '

$ws.Range("D20").Value = '0.0₃0897'
$ws.Range("E20").Value = '  -2.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.20%  '

$ws.Range("E25").Value = '  -1.42%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.26%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.01%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.31%  '

$ws.Range("E33").Value = '  +0.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0730'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.01%  '

$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.39%  '

$ws.Range("E39").Value = '  -1.68%  '

$ws.Range("E40").Value = '  -3.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.30'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").Value = '1.943.53'
$ws.Range("E43").Value = '  -2.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0278'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.81%  '

$ws.Range("E47").Value = '  -4.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '52.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("D49").Value = '2.478.53'
$ws.Range("E49").Value = '  -1.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.80%  '
